# Regenerate the localization-status report:
#  - the two outstanding files have moved from "Ready for handoff" to
#    "In Translation", so every Status cell that shows that value needs
#    updating (Overview!E:F and the per-locale Status column on each
#    language sheet).
#  - re-running the report also re-sizes the Status column now that the
#    cell text is shorter.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
# Best achievable ColumnWidth input for this host's column-width quantizer
# that lands closest to the report generator's newly-computed Status column
# width.
$newStatusColWidth = 12.5

# --- Overview sheet: Status is mirrored in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$usedOverview = $wsOverview.UsedRange
for ($r = 1; $r -le $usedOverview.Rows.Count; $r++) {
    foreach ($colLetter in @("E", "F")) {
        $cell = $wsOverview.Range($colLetter + $r)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

# --- Per-locale sheets: Status lives in column C ---
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        $cell = $ws.Range("C" + $r)
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = $newStatusColWidth
}
